$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet ("eggers"), so it lands
# at the end rather than Excel's default "insert before active sheet" spot.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "nr_studies"

# Header row (bold / centered, matching the header style used on the other sheets)
$ws.Range("A1").Value = "outcome"
$ws.Range("B1").Value = "n_effect_sizes"
$ws.Range("C1").Value = "k_studies"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").HorizontalAlignment = -4108  # xlCenter

$data = @(
    @("MR (tot)", 4, 3),
    @("CAIR", 26, 9),
    @("A/FF", 10, 5),
    @("GC/TA", 6, 4),
    @("PSC", 4, 3),
    @("CCK", 8, 3),
    @("Fun and enjoyment", 8, 5),
    @("ITC", 4, 4),
    @("SE", 5, 3),
    @("Cognitive skills (PYD)", 2, 2),
    @("Initiative (PYD)", 2, 2),
    @("Negative experiences (PYD)", 2, 2),
    @("Personal and social skills (PYD)", 2, 2),
    @("Goal setting (PYD)", 1, 1),
    @("MR (amotivation)", 4, 4),
    @("SP", 5, 2),
    @("MR (controlled)", 7, 5),
    @("MR (self-determined)", 11, 5),
    @("INSB", 6, 3),
    @("INTB", 6, 3),
    @("BPNF", 4, 2),
    @("BPNS", 4, 2),
    @("AGT/MC task-related", 3, 2),
    @("Parentsupp", 3, 2)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
